$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 25000
$ws.Range("H84").Value = 25000
$ws.Range("H132").Value = 3797.25
$ws.Range("I132").Value = 3729.8333
$ws.Range("J81").Value = 30000
$ws.Range("J84").Value = 30000
$ws.Range("K132").Value = 11189.4999
$ws.Range("L81").Value = 30000
$ws.Range("L84").Value = 90000
$ws.Range("M132").Value = -8659.499899999999
$ws.Range("N81").Value = -31996
$ws.Range("N84").Value = -99984

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 0
$ws.Range("H45").Value = 3581.2273
$ws.Range("H61").Value = 32260856
$ws.Range("H63").Value = 1749.5
$ws.Range("H66").Value = 1749.5
$ws.Range("H110").Value = 17545106
$ws.Range("H132").Value = 4617.343
$ws.Range("H136").Value = 32260856
$ws.Range("H141").Value = 78499
$ws.Range("I6").Value = 0
$ws.Range("I45").Value = 2207.9524
$ws.Range("I61").Value = 2518.3215
$ws.Range("I63").Value = 1749.5
$ws.Range("I66").Value = 1749.5
$ws.Range("I110").Value = 1302.8125
$ws.Range("I132").Value = 2211.9
$ws.Range("I136").Value = 2518.3215
$ws.Range("J6").Value = 0
$ws.Range("J141").Value = 78499
$ws.Range("K6").Value = 0
$ws.Range("K45").Value = 2207.9524
$ws.Range("K61").Value = 2518.3215
$ws.Range("K63").Value = 1749.5
$ws.Range("K66").Value = 8747.5
$ws.Range("K110").Value = 1302.8125
$ws.Range("K132").Value = 6635.700000000001
$ws.Range("K136").Value = 7554.9645
$ws.Range("L6").Value = 0
$ws.Range("L141").Value = 78499
$ws.Range("M45").Value = -1830.9524
$ws.Range("M61").Value = -2306.3215
$ws.Range("M63").Value = -1063.5
$ws.Range("M66").Value = -5315.5
$ws.Range("M110").Value = 742.1875
$ws.Range("M132").Value = -4105.700000000001
$ws.Range("M136").Value = -5004.9645
$ws.Range("N141").Value = -88859
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 9259513
$ws.Range("H24").Value = 4000
$ws.Range("H68").Value = 32499.5
$ws.Range("H71").Value = 32499.5
$ws.Range("H80").Value = 22727594
$ws.Range("H82").Value = 27856
$ws.Range("H83").Value = 22727594
$ws.Range("H85").Value = 27856
$ws.Range("H126").Value = 50052
$ws.Range("H134").Value = 4314363.5
$ws.Range("I22").Value = 11111346
$ws.Range("I24").Value = 2000
$ws.Range("I68").Value = 21999
$ws.Range("I71").Value = 21999
$ws.Range("I80").Value = 50000384
$ws.Range("I83").Value = 50000384
$ws.Range("I134").Value = 6100244
$ws.Range("J22").Value = 350
$ws.Range("J24").Value = 6000
$ws.Range("J80").Value = 269.75
$ws.Range("J82").Value = 55312
$ws.Range("J83").Value = 269.75
$ws.Range("J85").Value = 55312
$ws.Range("J126").Value = 50052
$ws.Range("K22").Value = 11111346
$ws.Range("K24").Value = 2000
$ws.Range("K68").Value = 21999
$ws.Range("K71").Value = 65997
$ws.Range("K80").Value = 50000384
$ws.Range("K83").Value = 250001920
$ws.Range("K134").Value = 18300732
$ws.Range("L22").Value = 350
$ws.Range("L24").Value = 6000
$ws.Range("L80").Value = 269.75
$ws.Range("L82").Value = 55312
$ws.Range("L83").Value = 1348.75
$ws.Range("L85").Value = 55312
$ws.Range("L126").Value = 50052
$ws.Range("M22").Value = -11111173
$ws.Range("M24").Value = -1765
$ws.Range("M68").Value = -21188
$ws.Range("M71").Value = -61941
$ws.Range("M80").Value = -49999386
$ws.Range("M83").Value = -249996928
$ws.Range("M134").Value = -18298197
$ws.Range("N22").Value = -696
$ws.Range("N24").Value = -6470
$ws.Range("N80").Value = -2265.75
$ws.Range("N82").Value = -56078
$ws.Range("N83").Value = -11332.75
$ws.Range("N85").Value = -57964
$ws.Range("N126").Value = -59932

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 524
$ws.Range("H134").Value = 4718.654
$ws.Range("I22").Value = 598.6667
$ws.Range("J22").Value = 300
$ws.Range("J134").Value = 5952
$ws.Range("K22").Value = 598.6667
$ws.Range("L22").Value = 300
$ws.Range("L134").Value = 17856
$ws.Range("M22").Value = -248.6667
$ws.Range("N22").Value = -1000
$ws.Range("N134").Value = -22926

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 125.166664
$ws.Range("H33").Value = 27778104
$ws.Range("H92").Value = 15386464
$ws.Range("H98").Value = 1502.9
$ws.Range("I7").Value = 110.2
$ws.Range("I33").Value = 66666748
$ws.Range("I92").Value = 1899
$ws.Range("J33").Value = 501.14285
$ws.Range("J92").Value = 19232606
$ws.Range("J98").Value = 1709.625
$ws.Range("K7").Value = 330.6
$ws.Range("K33").Value = 400000488
$ws.Range("K92").Value = 5697
$ws.Range("L33").Value = 3006.8571
$ws.Range("L92").Value = 57697818
$ws.Range("L98").Value = 5128.875
$ws.Range("M7").Value = -218.6
$ws.Range("M33").Value = -400000205
$ws.Range("M92").Value = -4449
$ws.Range("N33").Value = -3572.8571
$ws.Range("N92").Value = -57700314
$ws.Range("N98").Value = -8124.875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 63749.75
$ws.Range("H107").Value = 534048.0600000001
$ws.Range("H113").Value = 4830.2607
$ws.Range("H126").Value = 3551.6
$ws.Range("I52").Value = 37499.5
$ws.Range("I107").Value = 727474.75
$ws.Range("I113").Value = 3302.4194
$ws.Range("I126").Value = 1783.6923
$ws.Range("J107").Value = 2124.75
$ws.Range("J126").Value = 5466.8335
$ws.Range("K52").Value = 37499.5
$ws.Range("K107").Value = 727474.75
$ws.Range("K113").Value = 3302.4194
$ws.Range("K126").Value = 5351.0769
$ws.Range("L107").Value = 2124.75
$ws.Range("L126").Value = 16400.5005
$ws.Range("M52").Value = -37240.5
$ws.Range("M107").Value = -725554.75
$ws.Range("M113").Value = -1132.4194
$ws.Range("M126").Value = -2881.0769
$ws.Range("N107").Value = -5964.75
$ws.Range("N126").Value = -21340.5005

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3156.9443
$ws.Range("H105").Value = 37989.668
$ws.Range("H106").Value = 19349.625
$ws.Range("I40").Value = 2188.4
$ws.Range("J105").Value = 37989.668
$ws.Range("J106").Value = 19349.625
$ws.Range("K40").Value = 2188.4
$ws.Range("L105").Value = 37989.668
$ws.Range("L106").Value = 19349.625
$ws.Range("M40").Value = -2052.4
$ws.Range("N105").Value = -44977.668
$ws.Range("N106").Value = -21873.625

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2721.75
$ws.Range("H126").Value = 581.3684
$ws.Range("I96").Value = 3499.5
$ws.Range("I126").Value = 469.73334
$ws.Range("J96").Value = 1944
$ws.Range("K96").Value = 3499.5
$ws.Range("K126").Value = 1409.20002
$ws.Range("L96").Value = 1944
$ws.Range("M96").Value = -2126.5
$ws.Range("M126").Value = 1060.79998
$ws.Range("N96").Value = -4690
